# The sheet gains a new "effective amount" helper column at F: it resolves
# to E (if filled) else D, but only once both the category (old F) and
# sub-category (old G) drop-downs have a value. Inserting a column at F
# shifts the former F/G/H columns (category dropdown / sub-category dropdown
# / duplicate-category label) one slot to the right, to G/H/I - their
# values, data-validation dropdowns and existing cell formatting move with
# them automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Insert()

# New formula, shared down F1:F2, picking the effective transaction amount.
$ws.Range("F1:F2").Formula = '=if(And(G1<>"",H1<>""),if(E1<>"",E1,D1),)'

# F1:F2 pick up the existing highlighted ("selected") fill already used on
# this row (I1 still carries it), so copy that formatting over first...
$ws.Range("I1").Copy()
$ws.Range("F1:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ...then refine the numeric display: F1 mirrors the plain amount format
# used by E1/E2, F2 mirrors the signed-parentheses format used by D1/D2.
$ws.Range("F1").HorizontalAlignment = -4152
$ws.Range("F1").NumberFormat = "#,##0.00"
$ws.Range("F2").HorizontalAlignment = -4152
$ws.Range("F2").NumberFormat = "#,##0.00;(#,##0.00)"
